$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin price/volume data (GitHub Actions refresh)

$ws.Range("D2").Value = '26.078.60'
$ws.Range("E2").Value = '  +0.82%  '

$ws.Range("D3").Value = '1.597.92'
$ws.Range("E3").Value = '  +0.90%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.95'
$ws.Range("E5").Value = '  +0.94%  '

$ws.Range("E6").Value = '  -0.28%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.489'
$ws.Range("E7").Value = '  +2.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.248'
$ws.Range("E8").Value = '  +0.68%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0617'
$ws.Range("E9").Value = '  +0.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.34'
$ws.Range("E10").Value = '  +1.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0809'
$ws.Range("E11").Value = '  +2.42%  '

$ws.Range("D12").Value = '1.814.12'
$ws.Range("E12").Value = '  +0.64%  '

$ws.Range("D13").Value = '1.592.47'
$ws.Range("E13").Value = '  +0.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.04'
$ws.Range("E14").Value = '  +0.26%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.522'
$ws.Range("E15").Value = '  +2.37%  '

$ws.Range("D16").Value = '26.053.36'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.79'
$ws.Range("E17").Value = '  +3.43%  '

$ws.Range("D18").Value = '0.0₃0733'
$ws.Range("E18").Value = '  +1.04%  '

$ws.Range("E19").Value = '  -0.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '204.43'
$ws.Range("E20").Value = '  +6.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.28'
$ws.Range("E21").Value = '  +2.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.34'
$ws.Range("E22").Value = '  -0.23%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.06'
$ws.Range("E23").Value = '  +2.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("E24").Value = '  +16.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.13'
$ws.Range("E25").Value = '  +1.82%  '

$ws.Range("E26").Value = '  -0.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.123'
$ws.Range("E27").Value = '  -6.76%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.27'
$ws.Range("E28").Value = '  +1.31%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.57'
$ws.Range("E29").Value = '  +1.93%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0480'
$ws.Range("E30").Value = '  +1.94%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.17'
$ws.Range("E31").Value = '  +0.61%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.14'
$ws.Range("E32").Value = '  +0.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.97'
$ws.Range("E33").Value = '  -1.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.49'
$ws.Range("E34").Value = '  +0.31%  '

$ws.Range("E35").Value = '  -0.36%  '

$ws.Range("D36").Value = '1.133.98'
$ws.Range("E36").Value = '  +3.37%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0164'
$ws.Range("E37").Value = '  +8.95%  '

$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.806'
$ws.Range("E38").Value = '  +3.65%  '

$ws.Range("B39").Value = 'PaxDollar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  -0.14%  '

$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.503'
$ws.Range("E40").Value = '  +0.34%  '

$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.32'
$ws.Range("E41").Value = '  -1.20%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.778'
$ws.Range("E42").Value = '  -4.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.16'
$ws.Range("E43").Value = '  +0.43%  '

$ws.Range("D44").Value = '1.724.33'
$ws.Range("E44").Value = '  +0.51%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.57'
$ws.Range("E45").Value = '  -1.31%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.51'
$ws.Range("E46").Value = '  +0.31%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.23'
$ws.Range("E47").Value = '  +1.97%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0504'
$ws.Range("E48").Value = '  -0.83%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.408'
$ws.Range("E49").Value = '  +0.16%  '

$ws.Range("E50").Value = '  -0.22%  '

$ws.Range("D51").Value = '0.0₇0933'
$ws.Range("E51").Value = '  -16.60%  '

